$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.07759999999999
$ws.Range("C21").Value = -13.26190000000001
$ws.Range("C23").Value = -11.92960000000001
$ws.Range("C25").Value = -10.99299999999999
